# Generate Report for Handback
# Updates the localization-status workbook with the results of the latest
# handback run for a986b74c-6c90-4107-8854-0695d15b12ee.md:
#  - zh-cn sheet, row 7  -> handback received, but version mismatch error
#  - de-de sheet, row 7  -> handback received, but version mismatch error

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/682c5a9e52a9cba333307ef4a1132a9d2cfcfd18/e2e/a986b74c-6c90-4107-8854-0695d15b12ee.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bf382ee615e91c01cb295a2b026389ceeead537b/e2e/a986b74c-6c90-4107-8854-0695d15b12ee.md."
$targetDisplay = "a986b74c-6c90-4107-8854-0695d15b12ee.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/682c5a9e52a9cba333307ef4a1132a9d2cfcfd18/e2e/a986b74c-6c90-4107-8854-0695d15b12ee.md"

# ---------------- zh-cn sheet ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(7, 9).Value = $targetDisplay
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(7, 9), $targetUrl, "", "", $targetDisplay)
$wsZh.Cells.Item(7, 9).Style = "HyperLink"

$wsZh.Cells.Item(7, 10).Value = "a986b74c-6c90-4107-8854-0695d15b12ee.9759362da8f46991414a62462671b47d6edefa36.zh-cn.xlf"
$wsZh.Cells.Item(7, 11).Value = "2016-08-19 08:52:54"
$wsZh.Cells.Item(7, 16).Value = $errorDetail

# ---------------- de-de sheet ----------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(7, 9).Value = $targetDisplay
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(7, 9), $targetUrl, "", "", $targetDisplay)
$wsDe.Cells.Item(7, 9).Style = "HyperLink"

$wsDe.Cells.Item(7, 10).Value = "a986b74c-6c90-4107-8854-0695d15b12ee.9759362da8f46991414a62462671b47d6edefa36.de-de.xlf"
$wsDe.Cells.Item(7, 11).Value = "2016-08-19 08:53:02"
$wsDe.Cells.Item(7, 16).Value = $errorDetail
